$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions data refresh)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.279.96'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -3.32%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.303.90'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -5.24%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '545.35'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '171.98'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.73%  '
$ws.Range('E7').Value = '  -4.76%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.295.62'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -5.29%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.610'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.45%  '
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '52.55'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000264'
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '8.86'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.62%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.823.51'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -5.49%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '18.02'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.93%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.117'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.54%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.302.95'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -5.29%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.62'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -4.08%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '63.172.98'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.56%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.965'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '432.06'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.46'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +8.88%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.03'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '83.05'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.15%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '13.19'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +4.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.55'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.16%  '
$ws.Range('E28').Value = '  -4.71%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.60'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.84%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '29.08'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.88%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.33'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '11.30'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.87%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '571.13'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -6.50%  '
$ws.Range('E34').Value = '  -3.36%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '58.01'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.31%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.145'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('E38').Value = '  +3.92%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '34.92'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.68%  '
$ws.Range('E40').Value = '  -6.66%  '
$ws.Range('E41').Value = '  -4.51%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.107.99'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -7.43%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('E44').Value = '  -3.20%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.22'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0400'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.36%  '
$ws.Range('E47').Value = '  -4.27%  '
$ws.Range('E48').Value = '  -3.72%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.58'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -5.76%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '132.21'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.00%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.99'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.00%  '
